$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 605.1081
$ws.Range("I15").Value = 605.1081
$ws.Range("K15").Value = 1815.3243
$ws.Range("M15").Value = -1646.3243

$ws.Range("H33").Value = 2645969
$ws.Range("I33").Value = 4629875
$ws.Range("J33").Value = 761.5
$ws.Range("K33").Value = 4629875
$ws.Range("L33").Value = 761.5
$ws.Range("M33").Value = -4629646
$ws.Range("N33").Value = -1219.5

$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("M48").ClearContents()

$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("M56").ClearContents()

$ws.Range("H64").Value = 7796.4287
$ws.Range("I64").Value = 3986.5
$ws.Range("J64").Value = 8197.474
$ws.Range("K64").Value = 3986.5
$ws.Range("L64").Value = 8197.474
$ws.Range("M64").Value = -3738.5
$ws.Range("N64").Value = -8693.474

$ws.Range("H67").Value = 7796.4287
$ws.Range("I67").Value = 3986.5
$ws.Range("J67").Value = 8197.474
$ws.Range("K67").Value = 3986.5
$ws.Range("L67").Value = 8197.474
$ws.Range("M67").Value = -3128.5
$ws.Range("N67").Value = -9913.474

$ws.Range("H86").Value = 9570.25
$ws.Range("I86").Value = 5569
$ws.Range("J86").Value = 10141.857
$ws.Range("K86").Value = 5569
$ws.Range("L86").Value = 10141.857
$ws.Range("M86").Value = -4446
$ws.Range("N86").Value = -12387.857

$ws.Range("H89").Value = 9570.25
$ws.Range("I89").Value = 5569
$ws.Range("J89").Value = 10141.857
$ws.Range("K89").Value = 27845
$ws.Range("L89").Value = 50709.285
$ws.Range("M89").Value = -22229
$ws.Range("N89").Value = -61941.285

$ws.Range("H98").Value = 1887.9722
$ws.Range("I98").Value = 1726.303
$ws.Range("K98").Value = 1726.303
$ws.Range("M98").Value = -228.3030000000001

$ws.Range("H122").Value = 1887.9722
$ws.Range("I122").Value = 1726.303
$ws.Range("K122").Value = 5178.909000000001
$ws.Range("M122").Value = -2728.909000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1825196.9
$ws.Range("I2").Value = 2262689
$ws.Range("K2").Value = 2262689
$ws.Range("M2").Value = -2262576

$ws.Range("H21").Value = 5749.5
$ws.Range("I21").Value = 1699
$ws.Range("K21").Value = 1699
$ws.Range("M21").Value = -1325

$ws.Range("H45").Value = 5536438
$ws.Range("I45").Value = 8463930
$ws.Range("K45").Value = 8463930
$ws.Range("M45").Value = -8463553

$ws.Range("H61").Value = 3871.16
$ws.Range("I61").Value = 3757.0476
$ws.Range("J61").Value = 4470.25
$ws.Range("K61").Value = 3757.0476
$ws.Range("L61").Value = 4470.25
$ws.Range("M61").Value = -3545.0476
$ws.Range("N61").Value = -4894.25

$ws.Range("H74").Value = 105045.11
$ws.Range("I74").Value = 4903.6
$ws.Range("K74").Value = 4903.6
$ws.Range("M74").Value = -4029.6

$ws.Range("H77").Value = 105045.11
$ws.Range("I77").Value = 4903.6
$ws.Range("K77").Value = 24518
$ws.Range("M77").Value = -20150

$ws.Range("H102").Value = 3475092.5
$ws.Range("I102").Value = 4388321
$ws.Range("K102").Value = 4388321
$ws.Range("M102").Value = -4386699

$ws.Range("H116").Value = 1825196.9
$ws.Range("I116").Value = 2262689
$ws.Range("K116").Value = 2262689
$ws.Range("M116").Value = -2260395

$ws.Range("H132").Value = 2754.6453
$ws.Range("I132").Value = 1842.7391
$ws.Range("K132").Value = 5528.2173
$ws.Range("M132").Value = -2998.2173

$ws.Range("H136").Value = 3871.16
$ws.Range("I136").Value = 3757.0476
$ws.Range("J136").Value = 4470.25
$ws.Range("K136").Value = 11271.1428
$ws.Range("L136").Value = 13410.75
$ws.Range("M136").Value = -8721.1428
$ws.Range("N136").Value = -18510.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1825196.9
$ws.Range("I3").Value = 2262689
$ws.Range("K3").Value = 2262689
$ws.Range("M3").Value = -2262575

$ws.Range("H29").Value = 2472
$ws.Range("I29").Value = 680
$ws.Range("K29").Value = 680
$ws.Range("M29").Value = -391

$ws.Range("H31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()

$ws.Range("H33").Value = 30780.111
$ws.Range("I33").Value = 22004.2
$ws.Range("J33").Value = 41750
$ws.Range("K33").Value = 22004.2
$ws.Range("L33").Value = 41750
$ws.Range("M33").Value = -21668.2
$ws.Range("N33").Value = -42422

$ws.Range("H34").Value = 19995.715
$ws.Range("J34").Value = 19995.715
$ws.Range("L34").Value = 19995.715
$ws.Range("N34").Value = -20223.715

$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("M37").ClearContents()
$ws.Range("N37").ClearContents()

$ws.Range("H38").Value = 88428.57000000001
$ws.Range("I38").Value = 90000
$ws.Range("J38").Value = 88166.664
$ws.Range("K38").Value = 90000
$ws.Range("L38").Value = 88166.664
$ws.Range("N38").Value = -88998.664
$ws.Range("M38").Value = -89584

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 22646.78
$ws.Range("I31").Value = 2985.1667
$ws.Range("J31").Value = 28855.71
$ws.Range("K31").Value = 2985.1667
$ws.Range("L31").Value = 28855.71
$ws.Range("M31").Value = -2690.1667
$ws.Range("N31").Value = -29445.71

$ws.Range("H34").Value = 22646.78
$ws.Range("I34").Value = 2985.1667
$ws.Range("J34").Value = 28855.71
$ws.Range("K34").Value = 2985.1667
$ws.Range("L34").Value = 28855.71
$ws.Range("M34").Value = -2783.1667
$ws.Range("N34").Value = -29259.71

$ws.Range("H122").Value = 2994.2104
$ws.Range("I122").Value = 2994.2104
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 8982.6312
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -6532.6312
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 159896.83
$ws.Range("I132").Value = 144412.72
$ws.Range("K132").Value = 433238.16
$ws.Range("M132").Value = -430708.16

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 466.5
$ws.Range("I17").Value = 600
$ws.Range("J17").Value = 333
$ws.Range("K17").Value = 1800
$ws.Range("L17").Value = 999
$ws.Range("M17").Value = -1631
$ws.Range("N17").Value = -1337

$ws.Range("H63").Value = 10181.818
$ws.Range("I63").Value = 5000
$ws.Range("J63").Value = 10700
$ws.Range("K63").Value = 15000
$ws.Range("L63").Value = 32100
$ws.Range("N63").Value = -33598
$ws.Range("M63").Value = -14251

$ws.Range("H66").Value = 10181.818
$ws.Range("I66").Value = 5000
$ws.Range("J66").Value = 10700
$ws.Range("K66").Value = 45000
$ws.Range("L66").Value = 96300
$ws.Range("N66").Value = -103788
$ws.Range("M66").Value = -41256

$ws.Range("H107").Value = 1152
$ws.Range("I107").Value = 1000.625
$ws.Range("K107").Value = 3001.875
$ws.Range("M107").Value = -1081.875

$ws.Range("H118").Value = 1562.1666
$ws.Range("I118").Value = 1374.6
$ws.Range("K118").Value = 4123.799999999999
$ws.Range("M118").Value = -2880.799999999999

$ws.Range("H125").Value = 8687.5
$ws.Range("J125").Value = 9357.143
$ws.Range("L125").Value = 28071.429
$ws.Range("N125").Value = -37911.429

$ws.Range("H132").Value = 1976.8889
$ws.Range("I132").Value = 1904
$ws.Range("J132").Value = 1986
$ws.Range("K132").Value = 17136
$ws.Range("L132").Value = 17874
$ws.Range("M132").Value = -14606
$ws.Range("N132").Value = -22934

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 6502
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 6502
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 6502
$ws.Range("M10").ClearContents()
$ws.Range("N10").Value = -6840

$ws.Range("H80").Value = 16917708
$ws.Range("I80").Value = 21851204
$ws.Range("J80").Value = 2864.2856
$ws.Range("K80").Value = 21851204
$ws.Range("L80").Value = 2864.2856
$ws.Range("M80").Value = -21850206
$ws.Range("N80").Value = -4860.2856

$ws.Range("H83").Value = 16917708
$ws.Range("I83").Value = 21851204
$ws.Range("J83").Value = 2864.2856
$ws.Range("K83").Value = 109256020
$ws.Range("L83").Value = 14321.428
$ws.Range("M83").Value = -109251028
$ws.Range("N83").Value = -24305.428

$ws.Range("H108").Value = 75833.336
$ws.Range("J108").Value = 100000
$ws.Range("L108").Value = 100000
$ws.Range("N108").Value = -107680

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2182162.2
$ws.Range("I46").Value = 10873315
$ws.Range("J46").Value = 9373.9375
$ws.Range("K46").Value = 10873315
$ws.Range("L46").Value = 9373.9375
$ws.Range("M46").Value = -10873127
$ws.Range("N46").Value = -9749.9375

$ws.Range("H55").Value = 1728.1111
$ws.Range("I55").Value = 2868.4
$ws.Range("K55").Value = 2868.4
$ws.Range("M55").Value = -2695.4

$ws.Range("H132").Value = 7016.2383
$ws.Range("I132").Value = 7680.364
$ws.Range("K132").Value = 23041.092
$ws.Range("M132").Value = -20511.092

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 64885.6
$ws.Range("J46").Value = 69857
$ws.Range("L46").Value = 69857
$ws.Range("N46").Value = -70319

$ws.Range("H122").Value = 2432.1765
$ws.Range("I122").Value = 2290.3928
$ws.Range("J122").Value = 3093.8333
$ws.Range("K122").Value = 6871.178400000001
$ws.Range("L122").Value = 9281.499899999999
$ws.Range("M122").Value = -4421.178400000001
$ws.Range("N122").Value = -14181.4999

$ws.Range("H134").Value = 64885.6
$ws.Range("J134").Value = 69857
$ws.Range("L134").Value = 209571
$ws.Range("N134").Value = -214641
